$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-51: column B holds the "IPC RO" reading used to recompute the
# other three columns now that "IPC PO" (column C) is forced to 0:
#   C = 0
#   D = DELTA    = C - B = -B
#   E = DELTA^2  = D^2   = B^2
for ($r = 2; $r -le 51; $r++) {
    $b = $ws.Cells.Item($r, 2).Value2

    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = -$b
    $ws.Cells.Item($r, 5).Value = $b * $b
}

# Row 52 (TOTAL): C52 sums the DELTA column, E52 sums the DELTA^2 column.
$sumD = 0
$sumE = 0
for ($r = 2; $r -le 51; $r++) {
    $sumD += $ws.Cells.Item($r, 4).Value2
    $sumE += $ws.Cells.Item($r, 5).Value2
}
$ws.Cells.Item(52, 3).Value = $sumD
$ws.Cells.Item(52, 5).Value = $sumE

# Row 53 (MSE): average of the DELTA^2 column.
$ws.Cells.Item(53, 5).Value = $sumE / 50

Write-Output "done"
